$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Units sheet: single column A, rows 2-6 = 1..5
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Units")
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5
$ws.Range("A2:A6").RowHeight = 12.75
$ws.Range("A6").Select()

# ---------------------------------------------------------------------------
# Traits sheet: columns A & B, rows 2-5
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Traits")
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 1
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 2
$ws.Range("A2:B5").RowHeight = 12.75
$ws.Range("B6").Select()

# ---------------------------------------------------------------------------
# Soils sheet: column A, rows 2-5 (C2 already present/preserved)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Soils")
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("A2:A5").RowHeight = 12.75

# ---------------------------------------------------------------------------
# SoilLayers sheet: columns A & B, rows 2-5 (no explicit row height here)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("SoilLayers")
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 1
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 2
$ws.Range("B5").Select()

# ---------------------------------------------------------------------------
# Regions sheet: column A, rows 2-5
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Regions")
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("A2:A5").RowHeight = 12.75
$ws.Range("A5").Select()

# ---------------------------------------------------------------------------
# Sites sheet: columns A & B, rows 2-9
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Sites")
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 1
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 2
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 1
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 2
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 1
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = 2
$ws.Range("A2:B9").RowHeight = 12.75
$ws.Range("B10").Select()

# ---------------------------------------------------------------------------
# Fields sheet: columns A, B & C, rows 2-9
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Fields")
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = 4
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 5
$ws.Range("C6").Value = 1
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 6
$ws.Range("C7").Value = 2
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 7
$ws.Range("C8").Value = 3
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = 8
$ws.Range("C9").Value = 4
$ws.Range("A2:C9").RowHeight = 12.75
$ws.Range("C10").Select()

# ---------------------------------------------------------------------------
# Crops sheet: column A, rows 2-5
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Crops")
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("A2:A5").RowHeight = 12.75
$ws.Range("A5").Select()

# ---------------------------------------------------------------------------
# Researchers sheet: column A, rows 2-5
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Researchers")
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("A2:A5").RowHeight = 12.75
$ws.Range("A6").Select()

# ---------------------------------------------------------------------------
# MetStations sheet: column A, rows 2-5
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("MetStations")
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("A2:A5").RowHeight = 12.75
$ws.Range("A6").Select()

# ---------------------------------------------------------------------------
# Fertilizers sheet: column A, rows 2-5 - data populated now, but this sheet
# must end up as the active/selected tab, so its Activate()/Select() calls
# are issued at the very end of the script (after every other sheet's
# Select() calls) to ensure it remains the final active sheet.
# ---------------------------------------------------------------------------
$wsFertilizers = $wb.Worksheets.Item("Fertilizers")
$wsFertilizers.Range("A2").Value = 1
$wsFertilizers.Range("A3").Value = 2
$wsFertilizers.Range("A4").Value = 3
$wsFertilizers.Range("A5").Value = 4
$wsFertilizers.Range("A2:A5").RowHeight = 12.75

# ---------------------------------------------------------------------------
# Methods sheet: column A, rows 2-5
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Methods")
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("A2:A5").RowHeight = 12.75
$ws.Range("A6").Select()

# ---------------------------------------------------------------------------
# Now activate Fertilizers last so it ends up as the selected/active tab.
# ---------------------------------------------------------------------------
$wsFertilizers.Activate()
$wsFertilizers.Range("H5").Select()
